$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Gnai2"
$ws.Range("C2").Value = "Cnr1"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 162.399297
$ws.Range("H2").Value = 487.197891
$ws.Range("I2").Value = 0.3910371682630009
$ws.Range("J2").Value = 0.3910371682630009
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.810042
$ws.Range("N2").Value = 2.430126
$ws.Range("O2").Value = 0.1880856911129011
$ws.Range("P2").Value = 0.1880856911129011
$ws.Range("Q2").Value = 131.550251340474
$ws.Range("R2").Value = 1183.952262064266
$ws.Range("S2").Value = 0.07354849604357833
$ws.Range("T2").Value = 0.07354849604357833

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Gnai2"
$ws.Range("C3").Value = "Cnr1"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 162.399297
$ws.Range("H3").Value = 487.197891
$ws.Range("I3").Value = 0.3910371682630009
$ws.Range("J3").Value = 0.3910371682630009
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 3.496729
$ws.Range("N3").Value = 10.490187
$ws.Range("O3").Value = 0.811914308887099
$ws.Range("P3").Value = 0.8119143088870989
$ws.Range("Q3").Value = 567.8663313995129
$ws.Range("R3").Value = 5110.796982595616
$ws.Range("S3").Value = 0.3174886722194227
$ws.Range("T3").Value = 0.3174886722194226

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Gnai2"
$ws.Range("C4").Value = "Cnr1"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 65.41736466666667
$ws.Range("H4").Value = 196.252094
$ws.Range("I4").Value = 0.1575168212364948
$ws.Range("J4").Value = 0.1575168212364948
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.810042
$ws.Range("N4").Value = 2.430126
$ws.Range("O4").Value = 0.1880856911129011
$ws.Range("P4").Value = 0.1880856911129011
$ws.Range("Q4").Value = 52.99081290931601
$ws.Range("R4").Value = 476.917316183844
$ws.Range("S4").Value = 0.02962666018417342
$ws.Range("T4").Value = 0.02962666018417342

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Gnai2"
$ws.Range("C5").Value = "Cnr1"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 65.41736466666667
$ws.Range("H5").Value = 196.252094
$ws.Range("I5").Value = 0.1575168212364948
$ws.Range("J5").Value = 0.1575168212364948
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.496729
$ws.Range("N5").Value = 10.490187
$ws.Range("O5").Value = 0.811914308887099
$ws.Range("P5").Value = 0.8119143088870989
$ws.Range("Q5").Value = 228.7467961335087
$ws.Range("R5").Value = 2058.721165201578
$ws.Range("S5").Value = 0.1278901610523214
$ws.Range("T5").Value = 0.1278901610523214

# Row 6
$ws.Range("A6").Value = "M2"
$ws.Range("B6").Value = "Gnai2"
$ws.Range("C6").Value = "Cnr1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 126.3069433333333
$ws.Range("H6").Value = 378.92083
$ws.Range("I6").Value = 0.3041313008456065
$ws.Range("J6").Value = 0.3041313008456065
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.810042
$ws.Range("N6").Value = 2.430126
$ws.Range("O6").Value = 0.1880856911129011
$ws.Range("P6").Value = 0.1880856911129011
$ws.Range("Q6").Value = 102.31392899162
$ws.Range("R6").Value = 920.8253609245801
$ws.Range("S6").Value = 0.05720274590861153
$ws.Range("T6").Value = 0.05720274590861154

# Row 7
$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "Gnai2"
$ws.Range("C7").Value = "Cnr1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 126.3069433333333
$ws.Range("H7").Value = 378.92083
$ws.Range("I7").Value = 0.3041313008456065
$ws.Range("J7").Value = 0.3041313008456065
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.496729
$ws.Range("N7").Value = 10.490187
$ws.Range("O7").Value = 0.811914308887099
$ws.Range("P7").Value = 0.8119143088870989
$ws.Range("Q7").Value = 441.6611516550233
$ws.Range("R7").Value = 3974.95036489521
$ws.Range("S7").Value = 0.246928554936995
$ws.Range("T7").Value = 0.246928554936995

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Gnai2"
$ws.Range("C8").Value = "Cnr1"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 61.180387
$ws.Range("H8").Value = 183.541161
$ws.Range("I8").Value = 0.1473147096548978
$ws.Range("J8").Value = 0.1473147096548978
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.810042
$ws.Range("N8").Value = 2.430126
$ws.Range("O8").Value = 0.1880856911129011
$ws.Range("P8").Value = 0.1880856911129011
$ws.Range("Q8").Value = 49.558683046254
$ws.Range("R8").Value = 446.0281474162859
$ws.Range("S8").Value = 0.02770778897653781
$ws.Range("T8").Value = 0.02770778897653782

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Gnai2"
$ws.Range("C9").Value = "Cnr1"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 61.180387
$ws.Range("H9").Value = 183.541161
$ws.Range("I9").Value = 0.1473147096548978
$ws.Range("J9").Value = 0.1473147096548978
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 3.496729
$ws.Range("N9").Value = 10.490187
$ws.Range("O9").Value = 0.811914308887099
$ws.Range("P9").Value = 0.8119143088870989
$ws.Range("Q9").Value = 213.931233454123
$ws.Range("R9").Value = 1925.381101087107
$ws.Range("S9").Value = 0.11960692067836
$ws.Range("T9").Value = 0.11960692067836
